$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to previously entered daily case counts (column C) ---
$ws.Range("C343").Value = 104
$ws.Range("C390").Value = 38
$ws.Range("C402").Value = 88
$ws.Range("C410").Value = 82
$ws.Range("C412").Value = 159
$ws.Range("C413").Value = 84

# --- New day of data: row 414 (Date 2021-04-14 / serial 44300) ---
$ws.Range("C414").Value = 19
$ws.Range("E414").Value = 10
$ws.Range("F414").Value = 10
$ws.Range("G414").Value = 33

# L414/M414 are formatted as Text ("@") so a direct .Value assignment would
# store a text "0" instead of a real number. Round-trip the number format so
# the underlying stored value is numeric (matching how the rest of the
# column already holds real numbers), then restore the original Text format
# so the cell style is unchanged.
$ws.Range("L414").NumberFormat = "0"
$ws.Range("M414").NumberFormat = "0"
$ws.Range("L414").Value = 0
$ws.Range("M414").Value = 0
$ws.Range("L414").NumberFormat = "@"
$ws.Range("M414").NumberFormat = "@"
